$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.340.28'
$ws.Range("E2").Value = '  +1.78%  '

$ws.Range("D3").Value = '2.354.39'
$ws.Range("E3").Value = '  +2.61%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.54'
$ws.Range("E5").Value = '  +0.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.86'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  +0.19%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  +0.97%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.38'
$ws.Range("E10").Value = '  -0.60%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0787'
$ws.Range("E11").Value = '  +0.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.57'
$ws.Range("E12").Value = '  -3.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.122'
$ws.Range("E13").Value = '  +2.62%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.78'
$ws.Range("E14").Value = '  +0.05%  '

$ws.Range("D15").Value = '2.706.99'
$ws.Range("E15").Value = '  +1.99%  '

$ws.Range("D16").Value = '2.343.20'
$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.799'
$ws.Range("E17").Value = '  +2.35%  '

$ws.Range("D18").Value = '43.245.70'
$ws.Range("E18").Value = '  +1.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.21'
$ws.Range("E19").Value = '  -0.28%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.23'
$ws.Range("E20").Value = '  +3.74%  '

$ws.Range("D21").Value = '0.0₃0892'
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.21'
$ws.Range("E22").Value = '  +1.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.42'
$ws.Range("E23").Value = '  +0.50%  '

$ws.Range("E24").Value = '  -0.83%  '

$ws.Range("E25").Value = '  +0.09%  '

$ws.Range("E26").Value = '  +0.53%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.73'
$ws.Range("E27").Value = '  +0.74%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.24'
$ws.Range("E28").Value = '  +9.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.21'
$ws.Range("E29").Value = '  +2.20%  '

$ws.Range("E30").Value = '  -1.99%  '

$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.03'
$ws.Range("E32").Value = '  +1.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0732'
$ws.Range("E33").Value = '  +5.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.38'
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.84'
$ws.Range("E35").Value = '  +5.69%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.38'
$ws.Range("E36").Value = '  +0.19%  '

$ws.Range("E37").Value = '  -0.75%  '

$ws.Range("E38").Value = '  +1.13%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.79'
$ws.Range("E39").Value = '  +2.06%  '

$ws.Range("B40").Value = 'EnergySwap'
$ws.Range("C40").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.50'
$ws.Range("E40").Value = '  +15.84%  '

$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.45'
$ws.Range("E42").Value = '  -29.25%  '

$ws.Range("D43").Value = '1.943.28'
$ws.Range("E43").Value = '  -0.78%  '

$ws.Range("E44").Value = '  +1.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.98'
$ws.Range("E45").Value = '  -4.69%  '

$ws.Range("E46").Value = '  +2.04%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.74'
$ws.Range("E47").Value = '  -0.24%  '

$ws.Range("D48").Value = '2.572.91'
$ws.Range("E48").Value = '  +2.00%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '53.28'
$ws.Range("E49").Value = '  +0.67%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.79'
$ws.Range("E50").Value = '  -2.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.32'
$ws.Range("E51").Value = '  +0.67%  '
